$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TP_correlation_matrix")

# Update corrected July TP row (row 4) values
$ws.Range("B4").Value = 0.05
$ws.Range("C4").Value = 0.39
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 0.32
$ws.Range("F4").Value = 0.46
$ws.Range("G4").Value = 0.46

# Update the active view/selection state
$ws.Activate()
$ws.Range("I10").Select()
